$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1483
$ws.Range("D2").Value = 1483
$ws.Range("E2").Value = 1483
$ws.Range("F2").Value = 1483
$ws.Range("G2").Value = 1483
$ws.Range("H2").Value = 1483
$ws.Range("I2").Value = 1483
$ws.Range("J2").Value = 1483
$ws.Range("K2").Value = 1483
$ws.Range("L2").Value = 1483
$ws.Range("M2").Value = 1483
$ws.Range("N2").Value = 1483
$ws.Range("O2").Value = 1483
$ws.Range("P2").Value = 1483
$ws.Range("Q2").Value = 1483
$ws.Range("R2").Value = 1483
$ws.Range("S2").Value = 1483
$ws.Range("T2").Value = 1483
$ws.Range("U2").Value = 1483
$ws.Range("V2").Value = 1483
$ws.Range("W2").Value = 1483
$ws.Range("X2").Value = 1483
$ws.Range("Y2").Value = 1483
$ws.Range("Z2").Value = 1483
$ws.Range("AA2").Value = 1483
$ws.Range("AB2").Value = 1483
$ws.Range("AC2").Value = 1483
$ws.Range("AD2").Value = 1483
$ws.Range("AE2").Value = 1483
$ws.Range("AF2").Value = 1483
$ws.Range("AG2").Value = 1483
$ws.Range("AH2").Value = 1483
$ws.Range("AI2").Value = 1483
$ws.Range("AJ2").Value = 1483
$ws.Range("AK2").Value = 1483
$ws.Range("AL2").Value = 1483
$ws.Range("AM2").Value = 1483
# Row 3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4
$ws.Range("G3").Value = 159
$ws.Range("H3").Value = 400
$ws.Range("I3").Value = 108
$ws.Range("J3").Value = 181
$ws.Range("K3").Value = 46
$ws.Range("L3").Value = 11
$ws.Range("M3").Value = 44
$ws.Range("N3").Value = 109
$ws.Range("O3").Value = 37
$ws.Range("P3").Value = 14
$ws.Range("Q3").Value = 90
$ws.Range("R3").Value = 159
$ws.Range("S3").Value = 16
$ws.Range("T3").Value = 21
$ws.Range("U3").Value = 12
$ws.Range("V3").Value = 12
$ws.Range("W3").Value = 25
$ws.Range("X3").Value = 533
$ws.Range("Y3").Value = 564
$ws.Range("Z3").Value = 480
$ws.Range("AA3").Value = 297
$ws.Range("AB3").Value = 444
$ws.Range("AC3").Value = 539
$ws.Range("AD3").Value = 344
$ws.Range("AE3").Value = 305
$ws.Range("AF3").Value = 532
$ws.Range("AG3").Value = 611
$ws.Range("AH3").Value = 269
$ws.Range("AI3").Value = 381
$ws.Range("AJ3").Value = 289
$ws.Range("AK3").Value = 361
$ws.Range("AL3").Value = 445
# Row 4
$ws.Range("C4").Value = 2016
$ws.Range("D4").Value = 1.092380310182063
$ws.Range("G4").Value = 47.50573162508429
$ws.Range("H4").Value = 111.6392447741065
$ws.Range("I4").Value = 14.6621712744437
$ws.Range("J4").Value = 28.50708024275118
$ws.Range("K4").Value = 5.565745111260958
$ws.Range("L4").Value = 0.588671611598112
$ws.Range("M4").Value = 3.782872555630479
$ws.Range("N4").Value = 13.98853674983142
$ws.Range("O4").Value = 1.710721510451787
$ws.Range("P4").Value = 0.6749831422791639
$ws.Range("Q4").Value = 10.17397167902899
$ws.Range("R4").Value = 26.28590694538098
$ws.Range("S4").Value = 0.6284558327714093
$ws.Range("T4").Value = 1.113283884018881
$ws.Range("U4").Value = 0.6911665542818611
$ws.Range("V4").Value = 0.8186109238031019
$ws.Range("W4").Value = 2.507754551584626
$ws.Range("X4").Value = 0.3778254854598996
$ws.Range("Y4").Value = 0.1292882898205744
$ws.Range("Z4").Value = 0.02410559533737742
$ws.Range("AA4").Value = 0.002443420874401912
$ws.Range("AB4").Value = 0.01359440458305515
$ws.Range("AC4").Value = 0.05783632124471515
$ws.Range("AD4").Value = 0.3191713036101222
$ws.Range("AE4").Value = 0.004077921723339313
$ws.Range("AF4").Value = 0.3564420767356366
$ws.Range("AG4").Value = 0.2337153418023637
$ws.Range("AH4").Value = 0.001977483714266949
$ws.Range("AI4").Value = 0.005052931276727442
$ws.Range("AJ4").Value = 0.3338294163315963
$ws.Range("AK4").Value = 0.003499109577459762
$ws.Range("AL4").Value = 0.01322318488256458
$ws.Range("AM4").Value = 3396605.636547539
# Row 5
$ws.Range("C5").Value = 2016
$ws.Range("G5").Value = 31
$ws.Range("H5").Value = 11
$ws.Range("X5").Value = 0.1352657004830918
$ws.Range("Y5").Value = 0.1219512195121951
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0.008976660682226212
$ws.Range("AF5").Value = 0.09565217391304348
$ws.Range("AG5").Value = 0.2042253521126761
$ws.Range("AJ5").Value = 0.00816326530612245
# Row 6
$ws.Range("C6").Value = 1483
$ws.Range("D6").Value = 1353
$ws.Range("G6").Value = 17
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 95
$ws.Range("J6").Value = 83
$ws.Range("K6").Value = 831
$ws.Range("L6").Value = 1141
$ws.Range("M6").Value = 941
$ws.Range("N6").Value = 765
$ws.Range("O6").Value = 1087
$ws.Range("P6").Value = 1136
$ws.Range("Q6").Value = 771
$ws.Range("R6").Value = 27
$ws.Range("S6").Value = 1181
$ws.Range("T6").Value = 1038
$ws.Range("U6").Value = 1083
$ws.Range("V6").Value = 1059
$ws.Range("W6").Value = 886
$ws.Range("X6").Value = 1
$ws.Range("Y6").Value = 3
$ws.Range("Z6").Value = 831
$ws.Range("AA6").Value = 1141
$ws.Range("AB6").Value = 941
$ws.Range("AC6").Value = 765
$ws.Range("AD6").Value = 1
$ws.Range("AE6").Value = 1136
$ws.Range("AF6").Value = 1
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 1181
$ws.Range("AI6").Value = 1038
$ws.Range("AJ6").Value = 1
$ws.Range("AK6").Value = 1059
$ws.Range("AL6").Value = 886
$ws.Range("AM6").Value = 584
# Row 7
$ws.Range("C7").Value = 2016
$ws.Range("G7").Value = 3
$ws.Range("AC7").Value = 0
# Row 8
$ws.Range("C8").Value = 1483
$ws.Range("D8").Value = 1353
$ws.Range("G8").Value = 36
$ws.Range("H8").Value = 460
$ws.Range("I8").Value = 722
$ws.Range("J8").Value = 669
$ws.Range("K8").Value = 831
$ws.Range("L8").Value = 1141
$ws.Range("M8").Value = 941
$ws.Range("N8").Value = 765
$ws.Range("O8").Value = 1087
$ws.Range("P8").Value = 1136
$ws.Range("Q8").Value = 771
$ws.Range("R8").Value = 516
$ws.Range("S8").Value = 1181
$ws.Range("T8").Value = 1038
$ws.Range("U8").Value = 1083
$ws.Range("V8").Value = 1059
$ws.Range("W8").Value = 886
$ws.Range("X8").Value = 466
$ws.Range("Y8").Value = 669
$ws.Range("Z8").Value = 831
$ws.Range("AA8").Value = 1141
$ws.Range("AB8").Value = 941
$ws.Range("AC8").Value = 765
$ws.Range("AD8").Value = 627
$ws.Range("AE8").Value = 1136
$ws.Range("AF8").Value = 462
$ws.Range("AG8").Value = 516
$ws.Range("AH8").Value = 1181
$ws.Range("AI8").Value = 1038
$ws.Range("AJ8").Value = 630
$ws.Range("AK8").Value = 1059
$ws.Range("AL8").Value = 886
$ws.Range("AM8").Value = 584
# Row 9
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.3077320038129833
$ws.Range("G9").Value = 45.32907642513199
$ws.Range("H9").Value = 177.9291053453142
$ws.Range("I9").Value = 25.80268411024679
$ws.Range("J9").Value = 48.61975265205488
$ws.Range("K9").Value = 9.814760010817414
$ws.Range("L9").Value = 1.444681642891701
$ws.Range("M9").Value = 7.863978801666413
$ws.Range("N9").Value = 25.42524706553083
$ws.Range("O9").Value = 5.094312904896844
$ws.Range("P9").Value = 1.697688727030272
$ws.Range("Q9").Value = 18.66403578144393
$ws.Range("R9").Value = 40.84710948413129
$ws.Range("S9").Value = 1.857668056970902
$ws.Range("T9").Value = 2.498443281517798
$ws.Range("U9").Value = 1.565173064521483
$ws.Range("V9").Value = 1.69109087649545
$ws.Range("W9").Value = 4.600553040059128
$ws.Range("X9").Value = 0.4257468807024696
$ws.Range("Y9").Value = 0.1373271511276778
$ws.Range("Z9").Value = 0.04394085827308228
$ws.Range("AA9").Value = 0.007626517667080745
$ws.Range("AB9").Value = 0.02713584014058783
$ws.Range("AC9").Value = 0.07664531078664973
$ws.Range("AD9").Value = 0.4598706238918917
$ws.Range("AE9").Value = 0.02879855385072225
$ws.Range("AF9").Value = 0.4365026176660978
$ws.Range("AG9").Value = 0.2515131552269483
$ws.Range("AH9").Value = 0.005729023044261746
$ws.Range("AI9").Value = 0.01409937687392615
$ws.Range("AJ9").Value = 0.4561314359548418
$ws.Range("AK9").Value = 0.01096076583240338
$ws.Range("AL9").Value = 0.04275229409667387
$ws.Range("AM9").Value = 4738327.13204564
# Row 11
$ws.Range("C11").Value = 2016
$ws.Range("D11").Value = 4
$ws.Range("H11").Value = 672
$ws.Range("I11").Value = 123
$ws.Range("K11").Value = 48
$ws.Range("L11").Value = 11
$ws.Range("M11").Value = 47
$ws.Range("Q11").Value = 116
$ws.Range("S11").Value = 20
$ws.Range("T11").Value = 24
$ws.Range("U11").Value = 13
$ws.Range("AA11").Value = 0.1428571428571428
$ws.Range("AH11").Value = 0.09090909090909091
$ws.Range("AI11").Value = 0.25
$ws.Range("AK11").Value = 0.25
# Row 12
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 749
$ws.Range("H12").Value = 460
$ws.Range("I12").Value = 722
$ws.Range("J12").Value = 669
$ws.Range("K12").Value = 831
$ws.Range("L12").Value = 1141
$ws.Range("M12").Value = 941
$ws.Range("N12").Value = 765
$ws.Range("O12").Value = 1087
$ws.Range("P12").Value = 1136
$ws.Range("Q12").Value = 771
$ws.Range("R12").Value = 516
$ws.Range("S12").Value = 1181
$ws.Range("T12").Value = 1038
$ws.Range("U12").Value = 1083
$ws.Range("V12").Value = 1059
$ws.Range("W12").Value = 886
$ws.Range("X12").Value = 265
$ws.Range("Y12").Value = 669
$ws.Range("Z12").Value = 831
$ws.Range("AA12").Value = 1141
$ws.Range("AB12").Value = 941
$ws.Range("AC12").Value = 765
$ws.Range("AD12").Value = 627
$ws.Range("AE12").Value = 1136
$ws.Range("AF12").Value = 315
$ws.Range("AG12").Value = 516
$ws.Range("AH12").Value = 1181
$ws.Range("AI12").Value = 1038
$ws.Range("AJ12").Value = 630
$ws.Range("AK12").Value = 1059
$ws.Range("AL12").Value = 886
